$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 13:40"

# Update country name labels (column A) that changed rank/position in the shared string table
$ws.Range("A67").Value = "Nepal"
$ws.Range("A68").Value = "Marruecos"
$ws.Range("A109").Value = "Madagascar"
$ws.Range("A110").Value = "Guinea Ecuatorial"
$ws.Range("A111").Value = "Estonia"
$ws.Range("A112").Value = "Sudan del Sur"
$ws.Range("A201").Value = "Laos"
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Dominica"
$ws.Range("A204").Value = "Fiyi"
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A212").Value = "Montserrat"
$ws.Range("A213").Value = "Seychelles"

# Update numeric statistics (columns B,C,D,E,G,H) for rows whose data changed
$ws.Range("B4").Value = 2553686
$ws.Range("C4").Value = 730
$ws.Range("D4").Value = 1068768
$ws.Range("E4").Value = 1357269
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 127649
$ws.Range("B7").Value = 511478
$ws.Range("C7").Value = 2032
$ws.Range("D7").Value = 297013
$ws.Range("E7").Value = 198734
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 42
$ws.Range("H7").Value = 15731
$ws.Range("B13").Value = 220180
$ws.Range("C13").Value = 2456
$ws.Range("D13").Value = 180661
$ws.Range("E13").Value = 29155
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 125
$ws.Range("H13").Value = 10364
$ws.Range("B23").Value = 93663
$ws.Range("C23").Value = 879
$ws.Range("D23").Value = 77225
$ws.Range("E23").Value = 16328
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 110
$ws.Range("B29").Value = 61095
$ws.Range("C29").Value = 382
$ws.Range("D29").Value = 44126
$ws.Range("E29").Value = 16592
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 377
$ws.Range("B35").Value = 44391
$ws.Range("C35").Value = 688
$ws.Range("D35").Value = 34586
$ws.Range("E35").Value = 9461
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 344
$ws.Range("B43").Value = 31555
$ws.Range("C43").Value = 69
$ws.Range("D43").Value = 29000
$ws.Range("E43").Value = 593
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1962
$ws.Range("B48").Value = 26022
$ws.Range("C48").Value = 325
$ws.Range("D48").Value = 18530
$ws.Range("E48").Value = 5903
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 10
$ws.Range("H48").Value = 1589
$ws.Range("B50").Value = 24805
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 19137
$ws.Range("E50").Value = 5592
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 76
$ws.Range("B56").Value = 17580
$ws.Range("C56").Value = 58
$ws.Range("D56").Value = 16371
$ws.Range("E56").Value = 509
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 700
$ws.Range("B67").Value = 12309
$ws.Range("C67").Value = 554
$ws.Range("D67").Value = 2834
$ws.Range("E67").Value = 9447
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 28
$ws.Range("B68").Value = 11854
$ws.Range("C68").Value = 221
$ws.Range("D68").Value = 8700
$ws.Range("E68").Value = 2936
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 218
$ws.Range("B109").Value = 2005
$ws.Range("C109").Value = 83
$ws.Range("D109").Value = 907
$ws.Range("E109").Value = 1082
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 16
$ws.Range("B110").Value = 2001
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 515
$ws.Range("E110").Value = 1454
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 32
$ws.Range("B111").Value = 1986
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 1812
$ws.Range("E111").Value = 105
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 69
$ws.Range("B112").Value = 1942
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 224
$ws.Range("E112").Value = 1682
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 36
$ws.Range("B150").Value = 670
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 635
$ws.Range("E150").Value = 26
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 9
$ws.Range("B158").Value = 355
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 330
$ws.Range("E158").Value = 25
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0
$ws.Range("B163").Value = 255
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 102
$ws.Range("E163").Value = 144
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 9
$ws.Range("B186").Value = 75
$ws.Range("C186").Value = 5
$ws.Range("D186").Value = 38
$ws.Range("E186").Value = 37
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0
$ws.Range("B190").Value = 44
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 26
$ws.Range("E190").Value = 16
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 2
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 10
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 1
$ws.Range("B213").Value = 11
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 11
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0
